$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.053650188724249
$arrBF[0,2] = 1.050880821334307
$arrBF[0,3] = 1.059413228233384
$arrBF[0,4] = 1.068230967277937
$ws.Range("B2:F2").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.037365202626548
$arrIM[0,1] = 1.05866641924357
$arrIM[0,2] = 1.053633869327474
$arrIM[0,3] = 1.062142800196144
$arrIM[0,4] = 1.070936704493662
$ws.Range("I2:M2").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.055403096830602
$arrBF[0,2] = 1.052200904220148
$arrBF[0,3] = 1.061019939041431
$arrBF[0,4] = 1.070028090516585
$ws.Range("B3:F3").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.037738308848142
$arrIM[0,1] = 1.060066399894186
$arrIM[0,2] = 1.054765180773
$arrIM[0,3] = 1.063561731626831
$arrIM[0,4] = 1.072547329010809
$ws.Range("I3:M3").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.056534558412089
$arrBF[0,2] = 1.053052377789828
$arrBF[0,3] = 1.062057226333646
$arrBF[0,4] = 1.07118877299439
$ws.Range("B4:F4").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.03797720051891
$arrIM[0,1] = 1.060969202151512
$arrIM[0,2] = 1.055493970113599
$arrIM[0,3] = 1.064477033107041
$arrIM[0,4] = 1.073586889217516
$ws.Range("I4:M4").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.057009572302049
$arrBF[0,2] = 1.053409698958796
$arrBF[0,3] = 1.062492749013972
$arrBF[0,4] = 1.071676216978872
$ws.Range("B5:F5").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.038077027817218
$arrIM[0,1] = 1.06134801480646
$arrIM[0,2] = 1.055799585778427
$arrIM[0,3] = 1.064861157119962
$arrIM[0,4] = 1.074023305725721
$ws.Range("I5:M5").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.057089291382632
$arrBF[0,2] = 1.053469657614496
$arrBF[0,3] = 1.062565843129681
$arrBF[0,4] = 1.071758031649947
$ws.Range("B6:F6").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.038093754018946
$arrIM[0,1] = 1.061411576916443
$arrIM[0,2] = 1.055850855266772
$arrIM[0,3] = 1.064925614354194
$arrIM[0,4] = 1.07409654631287
$ws.Range("I6:M6").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.056540908114904
$arrBF[0,2] = 1.053057154827064
$arrBF[0,3] = 1.062063047964242
$arrBF[0,4] = 1.071195288217289
$ws.Range("B7:F7").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.037978536779534
$arrIM[0,1] = 1.060974266699098
$arrIM[0,2] = 1.055498056767319
$arrIM[0,3] = 1.064482168402531
$arrIM[0,4] = 1.073592723033375
$ws.Range("I7:M7").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.054243178211907
$arrBF[0,2] = 1.051327516386072
$arrBF[0,3] = 1.059956720352591
$arrBF[0,4] = 1.068838772663169
$ws.Range("B8:F8").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.037491822301842
$arrIM[0,1] = 1.059140193805012
$arrIM[0,2] = 1.054016878895823
$arrIM[0,3] = 1.0626229302272
$arrIM[0,4] = 1.071481573259737
$ws.Range("I8:M8").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.050172275822988
$arrBF[0,2] = 1.048258476013772
$arrBF[0,3] = 1.056226444579035
$arrBF[0,4] = 1.06466898643421
$ws.Range("B9:F9").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.036614604654652
$arrIM[0,1] = 1.055884216794039
$arrIM[0,2] = 1.051381576857985
$arrIM[0,3] = 1.059324422739516
$arrIM[0,4] = 1.067740793578947
$ws.Range("I9:M9").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.047442613310522
$arrBF[0,2] = 1.046197558995675
$arrBF[0,3] = 1.053726253064014
$arrBF[0,4] = 1.061876601015298
$ws.Range("B10:F10").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.036016406309811
$arrIM[0,1] = 1.053696618266694
$arrIM[0,2] = 1.049607114222979
$arrIM[0,3] = 1.057109690804467
$arrIM[0,4] = 1.065232223656361
$ws.Range("I10:M10").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.046256700794695
$arrBF[0,2] = 1.045301484091681
$arrBF[0,3] = 1.052640299790607
$arrBF[0,4] = 1.060664295946989
$ws.Range("B11:F11").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.035754152347265
$arrIM[0,1] = 1.05274517860502
$arrIM[0,2] = 1.048834444573802
$arrIM[0,3] = 1.056146791689279
$arrIM[0,4] = 1.064142310185643
$ws.Range("I11:M11").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.045815588513747
$arrBF[0,2] = 1.044968075995464
$arrBF[0,3] = 1.052236409042711
$arrBF[0,4] = 1.060213496227339
$ws.Range("B12:F12").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.035656249671516
$arrIM[0,1] = 1.052391126906149
$arrIM[0,2] = 1.048546780585648
$arrIM[0,3] = 1.05578852695394
$arrIM[0,4] = 1.063736898008372
$ws.Range("I12:M12").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.04591023656188
$arrBF[0,2] = 1.045039618977315
$arrBF[0,3] = 1.052323068731135
$arrBF[0,4] = 1.060310217006289
$ws.Range("B13:F13").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.035677272370015
$arrIM[0,1] = 1.052467101601693
$arrIM[0,2] = 1.048608515546876
$arrIM[0,3] = 1.055865403390602
$arrIM[0,4] = 1.063823886339172
$ws.Range("I13:M13").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.046220250896828
$arrBF[0,2] = 1.045273936092881
$arrBF[0,3] = 1.052606924740113
$arrBF[0,4] = 1.060627042913301
$ws.Range("B14:F14").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.035746069701631
$arrIM[0,1] = 1.052715925815533
$arrIM[0,2] = 1.04881067971239
$arrIM[0,3] = 1.056117189761061
$arrIM[0,4] = 1.064108810405183
$ws.Range("I14:M14").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.046411179378481
$arrBF[0,2] = 1.045418231228063
$arrBF[0,3] = 1.052781748510855
$arrBF[0,4] = 1.060822183612772
$ws.Range("B15:F15").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.035788392956566
$arrIM[0,1] = 1.052869148770915
$arrIM[0,2] = 1.048935151907832
$arrIM[0,3] = 1.056272243565719
$arrIM[0,4] = 1.064284285542453
$ws.Range("I15:M15").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.047521233648695
$arrBF[0,2] = 1.046256949852759
$arrBF[0,3] = 1.053798252227902
$arrBF[0,4] = 1.0619569891584
$ws.Range("B16:F16").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.036033742796249
$arrIM[0,1] = 1.053759672574339
$arrIM[0,2] = 1.049658301841867
$arrIM[0,3] = 1.05717351173939
$arrIM[0,4] = 1.065304478600858
$ws.Range("I16:M16").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.048216471001874
$arrBF[0,2] = 1.046782060759585
$arrBF[0,3] = 1.054434969365424
$arrBF[0,4] = 1.062667958150532
$ws.Range("B17:F17").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.036186776250516
$arrIM[0,1] = 1.054317141387843
$arrIM[0,2] = 1.050110751056782
$arrIM[0,3] = 1.057737797732539
$arrIM[0,4] = 1.065943420640554
$ws.Range("I17:M17").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.048621611236735
$arrBF[0,2] = 1.047087994441821
$arrBF[0,3] = 1.05480603315569
$arrBF[0,4] = 1.063082348328344
$ws.Range("B18:F18").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.0362757267273
$arrIM[0,1] = 1.054641899732523
$arrIM[0,2] = 1.050374241432264
$arrIM[0,3] = 1.058066560352872
$arrIM[0,4] = 1.066315750152157
$ws.Range("I18:M18").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.048759689683288
$arrBF[0,2] = 1.047192250145542
$arrBF[0,3] = 1.054932502075422
$arrBF[0,4] = 1.063223593287068
$ws.Range("B19:F19").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.036306003869
$arrIM[0,1] = 1.054752565879658
$arrIM[0,2] = 1.050464014700606
$arrIM[0,3] = 1.058178596591318
$arrIM[0,4] = 1.066442645233358
$ws.Range("I19:M19").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.04814191804378
$arrBF[0,2] = 1.04672575809946
$arrBF[0,3] = 1.054366689116301
$arrBF[0,4] = 1.06259170968293
$ws.Range("B20:F20").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.036170389446464
$arrIM[0,1] = 1.054257372128957
$arrIM[0,2] = 1.050062250609576
$arrIM[0,3] = 1.057677294165697
$arrIM[0,4] = 1.065874904996826
$ws.Range("I20:M20").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.04612897643023
$arrBF[0,2] = 1.045204951267176
$arrBF[0,3] = 1.052523350682127
$arrBF[0,4] = 1.06053375940247
$ws.Range("B21:F21").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.035725824162322
$arrIM[0,1] = 1.052642671198058
$arrIM[0,2] = 1.048751165720505
$arrIM[0,3] = 1.056043061660886
$arrIM[0,4] = 1.064024923227439
$ws.Range("I21:M21").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.044859812710289
$arrBF[0,2] = 1.044245478788727
$arrBF[0,3] = 1.051361358055323
$arrBF[0,4] = 1.059236968045566
$ws.Range("B22:F22").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.035443471731789
$arrIM[0,1] = 1.051623707668535
$arrIM[0,2] = 1.047923009711609
$arrIM[0,3] = 1.055012069630889
$arrIM[0,4] = 1.062858461608526
$ws.Range("I22:M22").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.045532962740023
$arrBF[0,2] = 1.044754428490688
$arrBF[0,3] = 1.051977643102262
$arrBF[0,4] = 1.059924700154893
$ws.Range("B23:F23").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.035593422519479
$arrIM[0,1] = 1.052164238581578
$arrIM[0,2] = 1.048362397381087
$arrIM[0,3] = 1.055558953206477
$arrIM[0,4] = 1.063477143609545
$ws.Range("I23:M23").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.048175606490225
$arrBF[0,2] = 1.046751199943748
$arrBF[0,3] = 1.054397543018277
$arrBF[0,4] = 1.06262616402871
$ws.Range("B24:F24").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.036177794899874
$arrIM[0,1] = 1.054284380531815
$arrIM[0,2] = 1.050084167158059
$arrIM[0,3] = 1.057704634284412
$arrIM[0,4] = 1.065905865362426
$ws.Range("I24:M24").Value = $arrIM

$arrBF = New-Object "object[,]" 1,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.051227404975001
$arrBF[0,2] = 1.049054473723363
$arrBF[0,3] = 1.057193102788797
$arrBF[0,4] = 1.06574911820426
$ws.Range("B25:F25").Value = $arrBF
$arrIM = New-Object "object[,]" 1,5
$arrIM[0,0] = 1.036843728702271
$arrIM[0,1] = 1.056728897247429
$arrIM[0,2] = 1.052065922405377
$arrIM[0,3] = 1.060179883316744
$arrIM[0,4] = 1.068710407826122
$ws.Range("I25:M25").Value = $arrIM

Write-Output "Updated vm_pu values for rows 2-25"
